$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 15) to the sheet, mirroring the existing rows' pattern.
$row = 15

# Copy the formatting (date style) from the cell above, then set the value.
$ws.Range("A" + ($row - 1)).Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 1).Value = 42620.888136574074

$ws.Cells.Item($row, 2).Value = 30
$ws.Cells.Item($row, 3).Value = 66
$ws.Cells.Item($row, 4).Value = 32
$ws.Cells.Item($row, 5).Value = 66
$ws.Cells.Item($row, 6).Value = 17
$ws.Cells.Item($row, 7).Value = 37773
$ws.Cells.Item($row, 8).Value = 19355
$ws.Cells.Item($row, 9).Value = 1122
$ws.Cells.Item($row, 10).Value = 223
$ws.Cells.Item($row, 11).Value = 109
$ws.Cells.Item($row, 12).Value = 14
$ws.Cells.Item($row, 13).Value = 3
$ws.Cells.Item($row, 14).Value = "Named"

$wb.Save()
